$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "RAFAEL"
$ws.Range("B4").Value = "PULIDO LOPEZ"
$ws.Range("C4").Value = "NO"

# D4 / E4 stay empty (blank text), matching the "NO" acompanante row above
# (row 3) which stores literal empty-string cells rather than leaving the
# cells absent. A bare "" assignment clears the cell entirely instead of
# writing an empty string, so force-text (leading apostrophe) then strip the
# resulting quote-prefix style to land on a plain empty text cell.
$ws.Range("D4").Value = "'"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'"
$ws.Range("E4").Style = "Normal"

$ws.Range("F4").Value = "21/9/2025, 19:13:29"
